# Add an "invalid_data1" worksheet, a duplicate of "data1" that contains
# one missing/invalid data point (D6 cleared), as used by the new invalid
# test cases added in test_read_data.py.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the data1 sheet and place the copy right after it.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "invalid_data1"

# Make the copy "invalid" by removing one of its data points.
$ws2.Range("D6").ClearContents()

# Restore / update the selection on the original sheet (no longer the
# active tab, and the whole data range is selected instead of F9).
$ws1.Range("A1:F8").Select()

# The new invalid_data1 sheet becomes the active tab/sheet, with cell
# H15 selected.
$ws2.Activate()
$ws2.Range("H15").Select()
